$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.874.15"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.57%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.519.60"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -2.31%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "587.13"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.79%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "171.22"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +3.32%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.44%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.520.35"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -2.27%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.137"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.12%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.70%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.11"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.08%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.340"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -5.16%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.47"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.19%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.991.30"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.85%  "
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.54%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "66.708.22"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.51%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.533.00"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.74%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.95"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +2.81%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.25"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.52%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "352.95"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.36%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.15"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.79%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.57"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.68%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +4.64%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "69.65"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.93%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.99"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.27%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.681.54"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.08%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.997"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.25%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0972"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.58%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "531.18"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.95%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.07"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.43%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.13%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.90%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.18%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.03%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.87%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "156.46"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.22%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.52"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.21%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.43"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +1.06%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.97%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.77"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.85%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.06"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.92%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.05%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.47"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +2.63%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "39.72"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.35%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "148.45"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.36%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.554"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -2.25%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0₆0275"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -4.48%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.66"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.71%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.68"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.43%  "
